$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-04 Tuesday" "2024-06-05 Wednesday"

Replace-Text "300×3=" "287×6="
Replace-Text "710×5=" "523×8="
Replace-Text "423×2=" "318×8="
Replace-Text "108×9=" "947×9="
Replace-Text "971×5=" "200×3="

Replace-Text "238×6=" "693×6="
Replace-Text "114×6=" "113×9="
Replace-Text "973×3=" "670×8="
Replace-Text "192×4=" "607×6="
Replace-Text "241×3=" "686×5="

Replace-Text "769×7=" "852×9="
Replace-Text "348×8=" "559×4="
Replace-Text "415×7=" "730×9="
Replace-Text "688×7=" "554×5="
Replace-Text "793×5=" "310×9="

Replace-Text "978×9=" "133×8="
Replace-Text "458×7=" "450×7="
Replace-Text "456×8=" "129×5="
Replace-Text "585×6=" "509×4="
Replace-Text "582×4=" "736×4="

Replace-Text "766×4=" "412×3="
Replace-Text "908×3=" "425×6="
Replace-Text "147×4=" "916×7="
Replace-Text "552×6=" "638×5="
Replace-Text "654×3=" "102×3="
